$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "57.310.40"
$ws.Range("E2").Value = "  -4.73%  "

# Row 3
$ws.Range("D3").Value = "3.140.90"
$ws.Range("E3").Value = "  -4.77%  "

# Row 4
$ws.Range("D4").Value = "0.997"
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").Value = "517.58"
$ws.Range("E5").Value = "  -6.98%  "

# Row 6
$ws.Range("D6").Value = "133.11"
$ws.Range("E6").Value = "  -6.84%  "

# Row 7
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.19%  "

# Row 8
$ws.Range("D8").Value = "3.137.00"
$ws.Range("E8").Value = "  -5.00%  "

# Row 9
$ws.Range("D9").Value = "0.445"
$ws.Range("E9").Value = "  -6.46%  "

# Row 10
$ws.Range("D10").Value = "7.18"
$ws.Range("E10").Value = "  -8.41%  "

# Row 11
$ws.Range("D11").Value = "0.107"
$ws.Range("E11").Value = "  -9.66%  "

# Row 12
$ws.Range("D12").Value = "0.380"
$ws.Range("E12").Value = "  -6.29%  "

# Row 13
$ws.Range("D13").Value = "3.672.10"
$ws.Range("E13").Value = "  -4.57%  "

# Row 14
$ws.Range("E14").Value = "  -2.09%  "

# Row 15
$ws.Range("D15").Value = "25.12"
$ws.Range("E15").Value = "  -7.17%  "

# Row 16
$ws.Range("D16").Value = "3.135.60"
$ws.Range("E16").Value = "  -4.93%  "

# Row 17
$ws.Range("D17").Value = "57.252.53"
$ws.Range("E17").Value = "  -4.49%  "

# Row 18
$ws.Range("D18").Value = "0.0000148"
$ws.Range("E18").Value = "  -10.46%  "

# Row 19
$ws.Range("D19").Value = "5.69"
$ws.Range("E19").Value = "  -6.87%  "

# Row 20
$ws.Range("D20").Value = "12.81"
$ws.Range("E20").Value = "  -9.76%  "

# Row 21
$ws.Range("D21").Value = "7.90"
$ws.Range("E21").Value = "  -7.44%  "

# Row 22
$ws.Range("D22").Value = "341.29"
$ws.Range("E22").Value = "  -8.75%  "

# Row 24
$ws.Range("D24").Value = "68.30"
$ws.Range("E24").Value = "  -7.18%  "

# Row 25
$ws.Range("D25").Value = "0.500"
$ws.Range("E25").Value = "  -7.89%  "

# Row 26
$ws.Range("D26").Value = "3.263.42"
$ws.Range("E26").Value = "  -5.20%  "

# Row 27
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.18%  "

# Row 28
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "0.0₃0927"
$ws.Range("E28").Value = "  -9.94%  "

# Row 29
$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").Value = "0.162"
$ws.Range("E29").Value = "  -6.07%  "

# Row 30
$ws.Range("D30").Value = "0.998"
$ws.Range("E30").Value = "  -0.09%  "

# Row 31
$ws.Range("D31").Value = "6.66"
$ws.Range("E31").Value = "  -6.97%  "

# Row 32
$ws.Range("E32").Value = "  -9.08%  "

# Row 33
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "6.81"
$ws.Range("E33").Value = "  -10.69%  "

# Row 34
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "21.46"
$ws.Range("E34").Value = "  -4.64%  "

# Row 35
$ws.Range("D35").Value = "1.18"
$ws.Range("E35").Value = "  -5.47%  "

# Row 36
$ws.Range("D36").Value = "4.80"
$ws.Range("E36").Value = "  -7.56%  "

# Row 37
$ws.Range("D37").Value = "157.32"
$ws.Range("E37").Value = "  -5.43%  "

# Row 38
$ws.Range("D38").Value = "6.14"
$ws.Range("E38").Value = "  -8.50%  "

# Row 39
$ws.Range("D39").Value = "1.37"
$ws.Range("E39").Value = "  -9.46%  "

# Row 40
$ws.Range("D40").Value = "25.51"
$ws.Range("E40").Value = "  -5.19%  "

# Row 41
$ws.Range("B41").Value = "RenzoRestakedETH"
$ws.Range("C41").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D41").Value = "3.163.71"
$ws.Range("E41").Value = "  -4.82%  "

# Row 42
$ws.Range("B42").Value = "Hedera"
$ws.Range("C42").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D42").Value = "0.0679"
$ws.Range("E42").Value = "  -8.20%  "

# Row 43
$ws.Range("D43").Value = "40.19"
$ws.Range("E43").Value = "  -3.85%  "

# Row 44
$ws.Range("D44").Value = "0.688"
$ws.Range("E44").Value = "  -8.15%  "

# Row 45
$ws.Range("B45").Value = "ONDO"
$ws.Range("C45").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D45").Value = "1.05"
$ws.Range("E45").Value = "  -5.79%  "

# Row 46
$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").Value = "3.84"
$ws.Range("E46").Value = "  -8.32%  "

# Row 47
$ws.Range("D47").Value = "0.995"
$ws.Range("E47").Value = "  -0.27%  "

# Row 48
$ws.Range("E48").Value = "  -8.54%  "

# Row 49
$ws.Range("D49").Value = "2.230.70"
$ws.Range("E49").Value = "  -5.04%  "

# Row 50
$ws.Range("D50").Value = "6.09"
$ws.Range("E50").Value = "  -6.57%  "

# Row 51
$ws.Range("D51").Value = "19.88"
$ws.Range("E51").Value = "  -6.25%  "
